$d = $word.ActiveDocument

# Collapse a range to the very end of the document's main story (this sits
# inside the existing trailing empty paragraph).
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd

# Build a WordprocessingML fragment for the new paragraphs that should be
# appended after the existing trailing empty paragraph. The fragment starts
# with an empty paragraph so that the pre-existing trailing empty paragraph
# at the insertion point is preserved (InsertXML merges the *last* paragraph
# of the inserted fragment into the paragraph at the insertion point, so a
# leading empty paragraph keeps that original paragraph looking untouched).
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData>' +
      '<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
        '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
      '</Relationships>' +
    '</pkg:xmlData></pkg:part>' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p/>' +
        '<w:p><w:r><w:t>Wednesday</w:t></w:r></w:p>' +
        '<w:p/>' +
        '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Viewgroups</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
        '<w:p/>' +
        '<w:p><w:r><w:t>Alignment is essential</w:t></w:r></w:p>' +
        '<w:p/>' +
        '<w:p><w:r><w:t xml:space="preserve">1 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LinearLayout</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
        '<w:p><w:r><w:t xml:space="preserve">2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RelativeLayout</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
        '<w:p><w:r><w:t xml:space="preserve">3 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ConstraintLayout</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
        '<w:p><w:r><w:t>4 Coordinatorlayout</w:t></w:r></w:p>' +
      '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part>' +
  '</pkg:package>'

$endRange.InsertXML($xml)
